# Lunes 13 de Setiembre
#
# The paragraph that used to read "Algo mas" is split into two runs
# ("Algo " / "mas") wrapped in a grammar-check proofErr pair, and a new
# paragraph "Otro mas" is added right after it (inheriting the
# "_GoBack" bookmark that used to close the original paragraph).

$d = $word.ActiveDocument

# --- Locate the paragraph that contains "Algo mas" -------------------
$found = $d.Content
$ok = $found.Find.Execute("Algo mas", $false, $false, $false, $false, `
                           $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find the text 'Algo mas' in the document"
}

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $found.Start -and $p.Range.End -ge $found.End) {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not resolve the paragraph containing 'Algo mas'"
}

# Whole paragraph (text + the trailing paragraph mark), so replacing it
# also carries along whatever sits right after the text (the _GoBack
# bookmark) into the freshly-inserted content.
$paraRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End)

# Preserve the paragraph's rsid so the split run keeps matching
# metadata instead of looking like a totally unrelated paragraph.
$paraXml = $targetPara.Range.WordOpenXML
$rsid = "00000000"
if ($paraXml -match 'w:rsidR="([0-9A-F]+)"') {
    $rsid = $matches[1]
}

$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$lang = '<w:rPr><w:lang w:val="es-ES"/></w:rPr>'

# Paragraph 1: "Algo mas" split into "Algo " + "mas" around a
# gramStart/gramEnd proofErr pair (as Word's grammar checker does when
# it flags a run that was just edited).
$p1 = '<w:p w:rsidR="' + $rsid + '" w:rsidRDefault="' + $rsid + '">' + `
        '<w:pPr>' + $lang + '</w:pPr>' + `
        '<w:r>' + $lang + '<w:t xml:space="preserve">Algo </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r>' + $lang + '<w:t>mas</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
      '</w:p>'

# Paragraph 2: brand-new "Otro mas" paragraph, ending with the
# bookmark pair that used to close the original paragraph.
$p2 = '<w:p>' + `
        '<w:pPr>' + $lang + '</w:pPr>' + `
        '<w:r>' + $lang + '<w:t>Otro mas</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
        '<w:bookmarkEnd w:id="0"/>' + `
      '</w:p>'

$body = $p1 + $p2

$packageXml = '<?xml version="1.0" standalone="yes"?>' + "`n" + `
  '<?mso-application progid="Word.Document"?>' + "`n" + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document ' + $W + '><w:body>' + $body + '</w:body></w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

# Replacing the whole paragraph range (not just the text) in one shot
# means the content that used to trail the text - the _GoBack bookmark
# - ends up trailing whatever we insert, i.e. inside the new "Otro mas"
# paragraph, exactly like Word would leave it after Enter + typing.
$paraRange.InsertXML($packageXml)

Write-Output "Split 'Algo mas' and added 'Otro mas' paragraph"
